$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "28.903.76"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  -1.29%  "

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "1.907.29"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  -1.26%  "

$ws.Cells.Item(4, 5).Value = "  +0.02%  "

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "324.73"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.25%  "

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +0.09%  "

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.4580"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -0.94%  "

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.3813"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -1.54%  "

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.07725"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -1.26%  "

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "0.9799"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +0.84%  "

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "22.07"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -2.30%  "

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "1.905.77"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -0.14%  "

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "5.677"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -1.64%  "

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "6.939"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -1.88%  "

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "0.07089"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +0.36%  "

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +0.04%  "

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "83.72"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -3.50%  "

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "0.000009454"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -3.28%  "

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "16.60"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -2.57%  "

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -0.05%  "

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "28.903.55"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -1.52%  "

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "5.317"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -2.82%  "

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "10.94"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -1.04%  "

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "2.097"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +0.16%  "

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "158.72"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +1.17%  "

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "19.00"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -2.08%  "

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "5.664"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -1.66%  "

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "117.61"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -0.65%  "

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "1.872"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +0.72%  "

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "0.09288"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -0.51%  "

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "0.8614"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -0.27%  "

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "5.086"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -1.70%  "

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "1.247"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -4.38%  "

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "3.064"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -0.43%  "

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "0.05706"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -1.14%  "

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "1.157"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +0.28%  "

$ws.Cells.Item(37, 5).Value = "  +0.06%  "

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "0.02040"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -2.00%  "

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "7.419"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -3.23%  "

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "0.5488"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -3.05%  "

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "0.1750"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -1.53%  "

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "2.888"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +6.28%  "

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "9.318"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -0.65%  "

$ws.Cells.Item(44, 2).Value = "RenderToken"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "2.130"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +2.51%  "

$ws.Cells.Item(45, 2).Value = "Decentraland"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "0.5163"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -2.10%  "

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "11.22"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -1.96%  "

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "0.06884"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +0.18%  "

$ws.Cells.Item(48, 2).Value = "PEPE"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "0.000002621"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -14.14%  "

$ws.Cells.Item(49, 2).Value = "NEARProtocol"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "1.775"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -1.86%  "

$ws.Cells.Item(50, 2).Value = "Quant"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "110.23"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -1.13%  "

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "0.2876"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -4.08%  "

Write-Output "Update complete"